$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws1.Cells.Item(2, 6).Value = 82
$ws1.Cells.Item(3, 6).Value = 209
$ws1.Cells.Item(4, 6).Value = 117
$ws1.Cells.Item(5, 6).Value = 1728
$ws1.Cells.Item(6, 6).Value = 3305
$ws1.Cells.Item(7, 6).Value = 982
$ws1.Cells.Item(8, 6).Value = 2149
$ws1.Cells.Item(9, 6).Value = 2067
$ws1.Cells.Item(10, 6).Value = 1081
$ws1.Cells.Item(11, 6).Value = 579
$ws1.Cells.Item(13, 6).Value = 1647
$ws1.Cells.Item(14, 6).Value = 371
$ws1.Cells.Item(18, 6).Value = 160
$ws1.Cells.Item(19, 6).Value = 1526
$ws1.Cells.Item(21, 6).Value = 679
$ws1.Cells.Item(22, 6).Value = 566
$ws1.Cells.Item(23, 6).Value = 12044
$ws1.Cells.Item(24, 6).Value = 12051
$ws1.Cells.Item(25, 6).Value = 887
$ws1.Cells.Item(29, 6).Value = 2
$ws1.Cells.Item(30, 6).Value = 309
$ws1.Cells.Item(31, 6).Value = 1888
$ws1.Cells.Item(33, 6).Value = 523

$ws3 = $wb.Worksheets.Item(3)  # 本地生活
$ws3.Cells.Item(2, 6).Value = 72

$ws4 = $wb.Worksheets.Item(4)  # 全部类型
$ws4.Cells.Item(2, 6).Value = 82
$ws4.Cells.Item(3, 6).Value = 72
$ws4.Cells.Item(4, 6).Value = 209
$ws4.Cells.Item(6, 6).Value = 117
$ws4.Cells.Item(7, 6).Value = 1728
$ws4.Cells.Item(8, 6).Value = 3305
$ws4.Cells.Item(9, 6).Value = 982
$ws4.Cells.Item(10, 6).Value = 2149
$ws4.Cells.Item(11, 6).Value = 2067
$ws4.Cells.Item(12, 6).Value = 1081
$ws4.Cells.Item(13, 6).Value = 579
$ws4.Cells.Item(15, 6).Value = 1647
$ws4.Cells.Item(16, 6).Value = 371
$ws4.Cells.Item(22, 6).Value = 160
$ws4.Cells.Item(23, 6).Value = 1526
$ws4.Cells.Item(25, 6).Value = 679
$ws4.Cells.Item(26, 6).Value = 566
$ws4.Cells.Item(27, 6).Value = 12044
$ws4.Cells.Item(28, 6).Value = 12051
$ws4.Cells.Item(29, 6).Value = 887
$ws4.Cells.Item(33, 6).Value = 2
$ws4.Cells.Item(34, 6).Value = 309
$ws4.Cells.Item(35, 6).Value = 1888
$ws4.Cells.Item(39, 6).Value = 523
